# Applies the "Added basic state pattern, main menu and name entry" edit:
#  - F. Complex intelligence...  paragraph recolored red -> orange (accent2)
#  - H. Image rotation... paragraph (3 runs) recolored red -> green (accent6)
#  - Optional-section mark total "5/12 marks" -> "7/12 marks"
#  - TODO list: remove "1. Make states polymorphic"
#  - TODO list: remove "Fix object layering (lower objects appear on top)"

$d = $word.ActiveDocument

function Find-ParagraphByText($doc, $needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        $p = $doc.Paragraphs.Item($i)
        if ($p.Range.Text.Contains($needle)) {
            return $p
        }
    }
    return $null
}

# ---------------------------------------------------------------------------
# 1. "F. Complex intelligence on an automated moving object (1 or 2 marks)"
#    red (FF0000) -> orange, Accent 2 (ED7D31 / themeColor="accent2")
# ---------------------------------------------------------------------------
$pF = Find-ParagraphByText $d "F. Complex intelligence on an automated moving object"
$xmlF = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p><w:pPr><w:rPr><w:color w:val="ED7D31" w:themeColor="accent2"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:color w:val="ED7D31" w:themeColor="accent2"/></w:rPr><w:t>F. Complex intelligence on an automated moving object (1 or 2 marks)</w:t></w:r></w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pF.Range.InsertXML($xmlF)

# ---------------------------------------------------------------------------
# 2. "H. Image rotation/manipulation using the CoordinateMapping object (1 mark)"
#    red (FF0000) -> green, Accent 6 (70AD47 / themeColor="accent6")
# ---------------------------------------------------------------------------
$pH = Find-ParagraphByText $d "H. Image rotation/manipulation using the"
$xmlH = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p><w:pPr><w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr></w:pPr>' + `
    '<w:r><w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr><w:t xml:space="preserve">H. Image rotation/manipulation using the </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr><w:t>CoordinateMapping</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:rPr><w:color w:val="70AD47" w:themeColor="accent6"/></w:rPr><w:t xml:space="preserve"> object (1 mark)</w:t></w:r></w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pH.Range.InsertXML($xmlH)

# ---------------------------------------------------------------------------
# 3. Optional-section total "5/12 marks" -> "7/12 marks"
# ---------------------------------------------------------------------------
$pMarks = Find-ParagraphByText $d "/12"
$xmlMarks = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + `
    '<w:p><w:r><w:t>7</w:t></w:r><w:r><w:t>/12</w:t></w:r><w:r><w:t xml:space="preserve"> marks</w:t></w:r></w:p>' + `
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$pMarks.Range.InsertXML($xmlMarks)

# ---------------------------------------------------------------------------
# 4. TODO list: drop "Fix object layering (lower objects appear on top)" and
#    "1. Make states polymorphic" entirely (delete from bottom up so the
#    other lookup isn't affected by shifting indices).
# ---------------------------------------------------------------------------
$pLayering = Find-ParagraphByText $d "Fix object layering"
$pLayering.Range.Delete()

$pStates = Find-ParagraphByText $d "Make states polymorphic"
$pStates.Range.Delete()

Write-Output "done"
